$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20..32 down from the row above (each row takes the old values
# of the row immediately above it), then set row 19 to the new values.
# Process bottom-up is not required since we read "D" etc via a captured
# snapshot below to avoid overwrite-before-read issues.

$rows = 19..32

# Capture current (pre-edit) values for D, J, K, L, M, P for rows 19..32
$colsToCopy = @("D", "J", "K", "L", "M", "P")
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($col in $colsToCopy) {
        $snapshot[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Rows 20..32 get the snapshot values of the row above (r-1)
for ($r = 32; $r -ge 20; $r--) {
    foreach ($col in $colsToCopy) {
        $ws.Range("$col$r").Value = $snapshot[$r - 1][$col]
    }
}

# Row 19 gets the new values
$ws.Range("D19").Value = 44680
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 13500
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13750
$ws.Range("P19").Value = 1058
